$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-16.59521376628768, 1.970115203258419, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 2.762031086568757, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, 2.519964540780618, -16.59521376628768, -16.59521376628768, 3.020295540505377, -16.59521376628768, 1.482024433424249, -16.59521376628768, 2.294792167850281, -16.59521376628768),
    @(-16.59521376628768, 0.9636329120315715, -16.59521376628768, -16.59521376628768, -16.59521376628768, 2.899885886586267, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(2.718636916842794, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, 1.855357168621295, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(3.746023796089267, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 1.186408487696678, -16.59521376628768, 1.882801356173054),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, 2.544250152848897, -16.59521376628768, 2.738904544602889, -16.59521376628768, -16.59521376628768, -16.59521376628768, 1.804118405152197),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, 1.999785047854507, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 1.914379560075965, 1.950191935501683),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 2.046037163309848),
    @(-16.59521376628768, -16.59521376628768, 4.321914251476782, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 2.192028826454424, -16.59521376628768),
    @(-16.59521376628768, 1.993065330978587, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 1.298354405270216, 1.878052012205086, 1.943160857141208, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 1.537993682267127, 1.101081413877504, 1.540300141634238, -16.59521376628768),
    @(-16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, -16.59521376628768, 1.529841448031934, 1.620589372225509, -16.59521376628768, -16.59521376628768),
    @(-16.59521376628768, 0.6677711076923921, -16.59521376628768, -16.59521376628768, 3.57124969941008, -16.59521376628768, 2.081862359666718, 1.031785677154957, -16.59521376628768, 2.270776191179874),
    @(-16.59521376628768, 1.499034057940374, -16.59521376628768, 2.711268523237472, -16.59521376628768, 2.551029220244141, 2.239396159124805, -16.59521376628768, -16.59521376628768, -16.59521376628768)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value2 = $rowVals[$j]
    }
}
